$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continental US message fix - update BASE MSRP (column D) and DPHF (column E) values

$ws.Range("D29").Value = 53100
$ws.Range("D30").Value = 55890
$ws.Range("D31").Value = 64365

$ws.Range("D32").Value = 86580
$ws.Range("E32").Value = 1025

$ws.Range("D33").Value = 91580
$ws.Range("E33").Value = 1025

$ws.Range("D34").NumberFormat = "#,##0"
$ws.Range("D34").Value = 99310
$ws.Range("E34").Value = 1025
